# New crime data collected - weekly CompStat refresh (28th Precinct, week of 5/15/2023-5/21/2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/number and reporting week dates (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Row 15 (Rape) ---
$ws.Range("M15").Value = -50

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 7.692307692307
$ws.Range("I16").Value = 65
$ws.Range("J16").Value = 63
$ws.Range("K16").Value = 3.174603174603
$ws.Range("L16").Value = 8.333333333333
$ws.Range("M16").Value = -17.721518987341
$ws.Range("N16").Value = -78.618421052631

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -6.25
$ws.Range("I17").Value = 96
$ws.Range("J17").Value = 86
$ws.Range("K17").Value = 11.627906976744
$ws.Range("L17").Value = 14.285714285714
$ws.Range("M17").Value = 37.142857142857
$ws.Range("N17").Value = -60.493827160493

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -68.421052631578
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 75
$ws.Range("K18").Value = -41.333333333333
$ws.Range("L18").Value = -20
$ws.Range("M18").Value = -8.333333333333
$ws.Range("N18").Value = -86.746987951807

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 134
$ws.Range("J19").Value = 131
$ws.Range("K19").Value = 2.290076335877
$ws.Range("L19").Value = 12.605042016806
$ws.Range("M19").Value = 21.818181818181
$ws.Range("N19").Value = -15.189873417721

# --- Row 20 (G.L.A.) : D20/E20 switch from numbers to the suppressed-stat text markers ---
$ws.Range("D15").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4104)

$ws.Range("E15").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4104)

$ws.Range("L20").Value = 0
$ws.Range("N20").Value = -67.213114754098

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -29.411764705882
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 364
$ws.Range("J21").Value = 377
$ws.Range("K21").Value = -3.448275862068
$ws.Range("L21").Value = 6.122448979591
$ws.Range("M21").Value = 13.75
$ws.Range("N21").Value = -67.557932263814

# --- Row 22 (Transit) ---
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = -56.25
$ws.Range("L22").Value = -12.5
$ws.Range("M22").Value = -12.5

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 33
$ws.Range("J23").Value = 36
$ws.Range("K23").Value = -8.333333333333
$ws.Range("L23").Value = 6.451612903225
$ws.Range("M23").Value = 200

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = 31.034482758620
$ws.Range("I24").Value = 481
$ws.Range("J24").Value = 427
$ws.Range("K24").Value = 12.646370023419
$ws.Range("L24").Value = 23.333333333333
$ws.Range("M24").Value = 32.872928176795

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 40
$ws.Range("I25").Value = 168
$ws.Range("J25").Value = 168
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5
$ws.Range("M25").Value = -8.695652173913

# --- Row 26 (UCR Rape*) ---
$ws.Range("I26").Value = 5
$ws.Range("K26").Value = 25
$ws.Range("L26").Value = -37.5

# --- Row 27 (Other Sex Crimes): C27 switches from the suppressed-stat text marker to a number ---
$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = -40

# --- Row 28 (Shooting Vic.) ---
$ws.Range("G28").Value = 1
$ws.Range("L28").Value = -25
$ws.Range("M28").Value = -33.333333333333
$ws.Range("N28").Value = -84.210526315789

# --- Row 29 (Shooting Inc.) ---
$ws.Range("G29").Value = 1
$ws.Range("L29").Value = -37.5
$ws.Range("M29").Value = -44.444444444444
$ws.Range("N29").Value = -84.848484848484
